$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.515.15'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.950.15'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '243.48'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '57.90'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('E10').Value = '  -7.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').Value = '2.234.39'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.824'
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '13.72'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.17'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '1.942.20'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').Value = '36.430.54'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.28'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '0.0₃0847'
$ws.Range('E20').Value = '  -4.42%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '228.23'
$ws.Range('E21').Value = '  -0.85%  '
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.47'
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.10'
$ws.Range('E26').Value = '  -2.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '160.40'
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.26'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.66'
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0608'
$ws.Range('E33').Value = '  -6.68%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.33'
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.24'
$ws.Range('E36').Value = '  +3.58%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  +10.75%  '
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.23'
$ws.Range('E39').Value = '  -15.43%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0972'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.16'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '1.360.61'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '87.29'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.10'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  -0.66%  '
$ws.Range('D50').Value = '2.126.03'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('E51').Value = '  -3.05%  '
